$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.300.96'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.646.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("E6").Value = '  +2.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.255'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0629'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.09'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0850'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.878.06'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.660.61'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.38%  '
$ws.Range("E14").Value = '  +1.11%  '
$ws.Range("E15").Value = '  +3.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.297.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0742'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '220.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.66%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.78'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0513'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.92%  '
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.58'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.300.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.31%  '
$ws.Range("E37").Value = '  +0.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.553'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.862'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.25%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.810'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.24'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.788.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.19'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.55%  '
$ws.Range("E47").Value = '  +2.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0106'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.92%  '
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.73'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0970'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.81%  '
